# Daily attendance processing - swap the order of names in the
# "Recorded By" (column G) cells that read "System, dnasr281@gmail.com"
# so they read "dnasr281@gmail.com, System" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
